$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.094.47"
$ws.Range("E2").Value = "  +0.46%  "

$ws.Range("D3").Value = "3.467.87"
$ws.Range("E3").Value = "  -0.73%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.52%  "

$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.607"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.97%  "

$ws.Range("D9").Value = "3.468.01"
$ws.Range("E9").Value = "  -0.66%  "

$ws.Range("E10").Value = "  +4.63%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.97"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.02%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.432"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.68%  "

$ws.Range("D13").Value = "4.072.25"
$ws.Range("E13").Value = "  -0.58%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.18%  "

$ws.Range("E15").Value = "  -0.14%  "

$ws.Range("D16").Value = "67.186.89"
$ws.Range("E16").Value = "  +0.64%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000176"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.89%  "

$ws.Range("D18").Value = "3.472.51"
$ws.Range("E18").Value = "  -0.64%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.26%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "387.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.76%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.89"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.97%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.998"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.12%  "

$ws.Range("E25").Value = "  +1.37%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.537"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.69%  "

$ws.Range("E27").Value = "  +1.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.77%  "

$ws.Range("E29").Value = "  -2.43%  "

$ws.Range("E30").Value = "  +0.44%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.56%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.41"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.19%  "

$ws.Range("E33").Value = "  +0.42%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.41"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.59%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.33"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.04%  "

$ws.Range("E36").Value = "  -0.02%  "

$ws.Range("E37").Value = "  +0.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "161.32"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.39%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.879"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.47%  "

$ws.Range("E40").Value = "  +10.66%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.83%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.81"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.00%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.59"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.52%  "

$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.803.08"
$ws.Range("E44").Value = "  -0.52%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "26.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.41%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "26.58"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.41%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0720"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.91%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "41.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.87%  "

$ws.Range("E49").Value = "  -0.20%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "331.41"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.88%  "

$ws.Range("E51").Value = "  -2.11%  "
